# Filter - Study - Test Suit
# The "Cases" row in the TabName table is renamed to "Participants"
# (CasesTab -> ParticipantsTab), and the active selection moves to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ParticipantsTab"

$ws.Range("A2").Select()
